# Daily auto-push edit: a new 3-hourly ranking sample for 2026/02/25 16:00
# was appended to the source log, landing (by date) right before the
# 2026/12/29 block that already exists in the sheet. Insert a new row at
# 882 so every following row shifts down by one, then fill in the new
# sample's four columns (date, weekday, hour, ranking).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 882:923 down to 883:924, leaving a blank row 882 behind.
$ws.Rows.Item(882).Insert()

# Write the new data row. The date/weekday columns are plain text in this
# sheet (not real dates), so force text entry with a leading apostrophe to
# stop Excel from auto-converting "2026/02/25" into a date serial, then
# clear the resulting quote-prefix formatting so the cell's style stays
# the plain default (matching every other data cell in the column).
$ws.Cells.Item(882, 1).Value = "'2026/02/25"
$ws.Cells.Item(882, 1).ClearFormats()
$ws.Cells.Item(882, 2).Value = "水"
$ws.Cells.Item(882, 3).Value = 16
$ws.Cells.Item(882, 4).Value = 201
